$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Update "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE,
#    and replace the full URL values with just the numeric match code.
# ----------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4023"

$batting.Range("D3").NumberFormat = "@"
$batting.Range("D3").Value = "4026"

# ----------------------------------------------------------------------
# 2. Update "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE,
#    and replace the full URL values with just the numeric match code.
# ----------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4023"

$bowling.Range("B3").NumberFormat = "@"
$bowling.Range("B3").Value = "4026"

# ----------------------------------------------------------------------
# 3. Insert a brand new "Player Info" sheet as the first sheet in the
#    workbook, with player biographical data.
# ----------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4685"
$playerInfo.Range("B2").Value = "Scott Christopher Kuggeleijn"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"
